# The single slide master's theme ("Integral") is recoloured to the
# stock "Office Theme" colour values (Design > Variants > Colors >
# "Office" in the UI). The font scheme and format scheme (fills, lines,
# effects) of the "Integral" and "Office Theme" themes are already
# byte-for-byte identical, so only the twelve colour-scheme slots need
# to change.
#
# PowerPoint's ThemeColor.RGB is a BGR-packed long (0x00BBGGRR), so each
# target RRGGBB hex value below is converted to B*65536 + G*256 + R.

function BGR([int]$r, [int]$g, [int]$b) { return $b * 65536 + $g * 256 + $r }

$p  = $ppt.ActivePresentation
$m  = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

$cs.Item(1).RGB  = (BGR 0x00 0x00 0x00)   # dk1      000000
$cs.Item(2).RGB  = (BGR 0xFF 0xFF 0xFF)   # lt1      FFFFFF
$cs.Item(3).RGB  = (BGR 0x44 0x54 0x6A)   # dk2      44546A
$cs.Item(4).RGB  = (BGR 0xE7 0xE6 0xE6)   # lt2      E7E6E6
$cs.Item(5).RGB  = (BGR 0x5B 0x9B 0xD5)   # accent1  5B9BD5
$cs.Item(6).RGB  = (BGR 0xED 0x7D 0x31)   # accent2  ED7D31
$cs.Item(7).RGB  = (BGR 0xA5 0xA5 0xA5)   # accent3  A5A5A5
$cs.Item(8).RGB  = (BGR 0xFF 0xC0 0x00)   # accent4  FFC000
$cs.Item(9).RGB  = (BGR 0x44 0x72 0xC4)   # accent5  4472C4
$cs.Item(10).RGB = (BGR 0x70 0xAD 0x47)   # accent6  70AD47
$cs.Item(11).RGB = (BGR 0x05 0x63 0xC1)   # hlink    0563C1
$cs.Item(12).RGB = (BGR 0x95 0x4F 0x72)   # folHlink 954F72
